# The commit adds a new weekly price-record row for "Ají" (Inferno / Primera /
# Región de Arica y Parinacota) into the data table, inserted right before the
# existing row 334. All rows from the old row 334 onward shift down by one
# row, and the sheet's used range grows from A1:R390 to A1:R391.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 334 (pushes old row 334.. down to 335..)
$ws.Rows.Item(334).Insert()

# Populate the newly inserted row 334 with the new record's values.
$ws.Cells.Item(334, 1).Value2  = 9
$ws.Cells.Item(334, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(334, 3).Value2  = "Metropolitana"
$ws.Cells.Item(334, 4).Value2  = 45015
$ws.Cells.Item(334, 5).Value2  = 13
$ws.Cells.Item(334, 6).Value2  = 100112021
$ws.Cells.Item(334, 7).Value2  = "Ají"
$ws.Cells.Item(334, 8).Value2  = "Inferno"
$ws.Cells.Item(334, 9).Value2  = "Primera"
$ws.Cells.Item(334, 10).Value2 = 70
$ws.Cells.Item(334, 11).Value2 = 13000
$ws.Cells.Item(334, 12).Value2 = 15000
$ws.Cells.Item(334, 13).Value2 = 14000
$ws.Cells.Item(334, 14).Value2 = "$/caja 10 kilos"
$ws.Cells.Item(334, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(334, 16).Value2 = 1400
$ws.Cells.Item(334, 17).Value2 = 10
$ws.Cells.Item(334, 18).Value2 = "Hortaliza"
